# Daily Update - refresh monitor listing data (10 keyword rotations)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("E3").Value = "'188670"

# Row 12
$ws.Range("B12").Value = '벤큐 GW2780 아이케어'
$ws.Range("C12").Value = 'https://search.shopping.naver.com/gate.nhn?id=11846746201'
$ws.Range("D12").Value = 'https://shopping-phinf.pstatic.net/main_1184674/11846746201.20211118104344.jpg'
$ws.Range("E12").Value = "'209000"

# Row 13
$ws.Range("B13").Value = 'LG전자 27TN600S'
$ws.Range("C13").Value = 'https://search.shopping.naver.com/gate.nhn?id=21745903830'
$ws.Range("D13").Value = 'https://shopping-phinf.pstatic.net/main_2174590/21745903830.20200824150453.jpg'
$ws.Range("E13").Value = "'341800"
$ws.Range("I13").Value = 'LG전자'
$ws.Range("J13").Value = 'LG전자'

# Row 14
$ws.Range("B14").Value = '벤큐 XL2411K'
$ws.Range("C14").Value = 'https://search.shopping.naver.com/gate.nhn?id=24196764522'
$ws.Range("D14").Value = 'https://shopping-phinf.pstatic.net/main_2419676/24196764522.20211117182333.jpg'
$ws.Range("E14").Value = "'279000"
$ws.Range("I14").Value = '벤큐'
$ws.Range("J14").Value = '벤큐'

# Row 17
$ws.Range("B17").Value = '삼성전자 스마트모니터 M7 S43AM700'
$ws.Range("C17").Value = 'https://search.shopping.naver.com/gate.nhn?id=27215825524'
$ws.Range("D17").Value = 'https://shopping-phinf.pstatic.net/main_2721582/27215825524.20210520173025.jpg'
$ws.Range("E17").Value = "'649000"

# Row 18
$ws.Range("B18").Value = 'ASUS VA24DQLB'
$ws.Range("C18").Value = 'https://search.shopping.naver.com/gate.nhn?id=23257141490'
$ws.Range("D18").Value = 'https://shopping-phinf.pstatic.net/main_2325714/23257141490.20200720175614.jpg'
$ws.Range("E18").Value = "'219000"
$ws.Range("I18").Value = 'ASUS'
$ws.Range("J18").Value = 'ASUS'

# Row 20
$ws.Range("B20").Value = '삼성전자 스마트모니터 M7 S32AM700'
$ws.Range("C20").Value = 'https://search.shopping.naver.com/gate.nhn?id=25524333522'
$ws.Range("D20").Value = 'https://shopping-phinf.pstatic.net/main_2552433/25524333522.20210203133321.jpg'
$ws.Range("E20").Value = "'499000"
$ws.Range("I20").Value = '스마트모니터'
$ws.Range("J20").Value = '삼성전자'

# Row 36
$ws.Range("B36").Value = '삼성전자 오디세이 G5 C32G54T'
$ws.Range("C36").Value = 'https://search.shopping.naver.com/gate.nhn?id=23896004523'
$ws.Range("D36").Value = 'https://shopping-phinf.pstatic.net/main_2389600/23896004523.20210203132926.jpg'
$ws.Range("E36").Value = "'420000"
$ws.Range("I36").Value = '오디세이'
$ws.Range("J36").Value = '삼성전자'

# Row 37
$ws.Range("B37").Value = '한성컴퓨터 TFG32Q07P 75'
$ws.Range("C37").Value = 'https://search.shopping.naver.com/gate.nhn?id=28655748554'
$ws.Range("D37").Value = 'https://shopping-phinf.pstatic.net/main_2865574/28655748554.20210831152013.jpg'
$ws.Range("E37").Value = "'259000"
$ws.Range("I37").Value = '한성컴퓨터'
$ws.Range("J37").Value = '한성컴퓨터'

# Row 61
$ws.Range("E61").Value = "'237000"

# Row 64
$ws.Range("B64").Value = '삼성전자 오디세이 G9 C49G95T'
$ws.Range("C64").Value = 'https://search.shopping.naver.com/gate.nhn?id=23255316490'
$ws.Range("D64").Value = 'https://shopping-phinf.pstatic.net/main_2325531/23255316490.20210203133605.jpg'
$ws.Range("E64").Value = "'1690000"
$ws.Range("I64").Value = '오디세이'
$ws.Range("J64").Value = '삼성전자'

# Row 65
$ws.Range("B65").Value = '벤큐 XL2546K'
$ws.Range("C65").Value = 'https://search.shopping.naver.com/gate.nhn?id=24235203522'
$ws.Range("D65").Value = 'https://shopping-phinf.pstatic.net/main_2423520/24235203522.20211117182410.jpg'
$ws.Range("E65").Value = "'649000"
$ws.Range("I65").Value = '벤큐'
$ws.Range("J65").Value = '벤큐'

# Row 67
$ws.Range("B67").Value = '알파스캔 AOC 27B2 보더리스 75 시력보호'
$ws.Range("C67").Value = 'https://search.shopping.naver.com/gate.nhn?id=21720504796'
$ws.Range("D67").Value = 'https://shopping-phinf.pstatic.net/main_2172050/21720504796.20210310171806.jpg'
$ws.Range("E67").Value = "'219000"
$ws.Range("I67").Value = '알파스캔'
$ws.Range("J67").Value = '알파스캔'

# Row 68
$ws.Range("B68").Value = '한성컴퓨터 TFG39Q14V 144'
$ws.Range("C68").Value = 'https://search.shopping.naver.com/gate.nhn?id=26826361522'
$ws.Range("D68").Value = 'https://shopping-phinf.pstatic.net/main_2682636/26826361522.20210419161946.jpg'
$ws.Range("E68").Value = "'499000"
$ws.Range("I68").Value = '한성컴퓨터'
$ws.Range("J68").Value = '한성컴퓨터'

# Row 69
$ws.Range("B69").Value = '벤큐 ZOWIE XL2731'
$ws.Range("C69").Value = 'https://search.shopping.naver.com/gate.nhn?id=22435628535'
$ws.Range("D69").Value = 'https://shopping-phinf.pstatic.net/main_2243562/22435628535.20211126161127.jpg'
$ws.Range("E69").Value = "'419000"
$ws.Range("I69").Value = '벤큐'
$ws.Range("J69").Value = '벤큐'

# Row 73
$ws.Range("E73").Value = "'567000"

# Row 91
$ws.Range("B91").Value = 'LG전자 울트라와이드 29WP500'
$ws.Range("C91").Value = 'https://search.shopping.naver.com/gate.nhn?id=26886077522'
$ws.Range("D91").Value = 'https://shopping-phinf.pstatic.net/main_2688607/26886077522.20210524134552.jpg'
$ws.Range("E91").Value = "'259000"
$ws.Range("I91").Value = '울트라와이드'
$ws.Range("J91").Value = 'LG전자'

# Row 92
$ws.Range("B92").Value = '한성컴퓨터 TFG27Q14P 144'
$ws.Range("C92").Value = 'https://search.shopping.naver.com/gate.nhn?id=27327723522'
$ws.Range("D92").Value = 'https://shopping-phinf.pstatic.net/main_2732772/27327723522.20210527095004.jpg'
$ws.Range("E92").Value = "'379000"
$ws.Range("I92").Value = '한성컴퓨터'
$ws.Range("J92").Value = '한성컴퓨터'

# Row 93
$ws.Range("B93").Value = '삼성전자 삼성 U32R590'
$ws.Range("C93").Value = 'https://search.shopping.naver.com/gate.nhn?id=17650306747'
$ws.Range("D93").Value = 'https://shopping-phinf.pstatic.net/main_1765030/17650306747.20210203134432.jpg'
$ws.Range("E93").Value = "'399000"
$ws.Range("I93").Value = '삼성'
$ws.Range("J93").Value = '삼성전자'

# Row 94
$ws.Range("B94").Value = '삼성전자 삼성 C27F391'
$ws.Range("C94").Value = 'https://search.shopping.naver.com/gate.nhn?id=9681100715'
$ws.Range("D94").Value = 'https://shopping-phinf.pstatic.net/main_9681100/9681100715.20200915114554.jpg'
$ws.Range("E94").Value = "'238000"
$ws.Range("I94").Value = '삼성'
$ws.Range("J94").Value = '삼성전자'

# Row 95
$ws.Range("B95").Value = 'DELL 울트라샤프 U2720Q'
$ws.Range("C95").Value = 'https://search.shopping.naver.com/gate.nhn?id=21752731630'
$ws.Range("D95").Value = 'https://shopping-phinf.pstatic.net/main_2175273/21752731630.20200327122054.jpg'
$ws.Range("E95").Value = "'778990"
$ws.Range("I95").Value = '울트라샤프'
$ws.Range("J95").Value = 'DELL'

# Row 100
$ws.Range("B100").Value = '삼성전자 삼성 C27F390'
$ws.Range("C100").Value = 'https://search.shopping.naver.com/gate.nhn?id=9489557554'
$ws.Range("D100").Value = 'https://shopping-phinf.pstatic.net/main_9489557/9489557554.20210203132811.jpg'
$ws.Range("E100").Value = "'238000"
$ws.Range("I100").Value = '삼성'
$ws.Range("J100").Value = '삼성전자'

# Row 101
$ws.Range("B101").Value = '벤큐 모비우스 EX2710S'
$ws.Range("C101").Value = 'https://search.shopping.naver.com/gate.nhn?id=27862189523'
$ws.Range("D101").Value = 'https://shopping-phinf.pstatic.net/main_2786218/27862189523.20211116103224.jpg'
$ws.Range("E101").Value = "'369000"
$ws.Range("I101").Value = '벤큐'
$ws.Range("J101").Value = '벤큐'

